$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.019999999999999
$ws.Range("C2").Value = 1.057518251446868
$ws.Range("D2").Value = 1.056715957242074
$ws.Range("E2").Value = 1.062893469647623
$ws.Range("F2").Value = 1.071854599782754
$ws.Range("I2").Value = 1.042107760283332
$ws.Range("J2").Value = 1.062514257117994
$ws.Range("K2").Value = 1.059452906078154
$ws.Range("L2").Value = 1.065613583238335
$ws.Range("M2").Value = 1.07455066508777
$ws.Range("N2").Value = 1.064023149235191

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.05911502416933
$ws.Range("D3").Value = 1.058142751221481
$ws.Range("E3").Value = 1.064356979464462
$ws.Range("F3").Value = 1.073502876126512
$ws.Range("I3").Value = 1.042479212575983
$ws.Range("J3").Value = 1.063760419306233
$ws.Range("K3").Value = 1.060691834741804
$ws.Range("L3").Value = 1.066890369012287
$ws.Range("M3").Value = 1.076013523524804
$ws.Range("N3").Value = 1.065271081116673

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.060145989020426
$ws.Range("D4").Value = 1.059064097709289
$ws.Range("E4").Value = 1.065302196934802
$ws.Range("F4").Value = 1.07456785146581
$ws.Range("I4").Value = 1.042716730368081
$ws.Range("J4").Value = 1.064564145330729
$ws.Range("K4").Value = 1.061491101628025
$ws.Range("L4").Value = 1.067714250061572
$ws.Range("M4").Value = 1.07695804147535
$ws.Range("N4").Value = 1.066075948524313

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.06057887737803
$ws.Range("D5").Value = 1.059450989022539
$ws.Range("E5").Value = 1.065699151008578
$ws.Range("F5").Value = 1.075015201179158
$ws.Range("I5").Value = 1.04281590657633
$ws.Range("J5").Value = 1.064901411444186
$ws.Range("K5").Value = 1.061826545971066
$ws.Range("L5").Value = 1.068060070709025
$ws.Range("M5").Value = 1.077354635239489
$ws.Range("N5").Value = 1.066413693594336

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.060651530499894
$ws.Range("D6").Value = 1.059515924057079
$ws.Range("E6").Value = 1.065765777239053
$ws.Range("F6").Value = 1.07509029199197
$ws.Range("I6").Value = 1.042832519122917
$ws.Range("J6").Value = 1.064958003763005
$ws.Range("K6").Value = 1.061882835501232
$ws.Range("L6").Value = 1.068118104145238
$ws.Range("M6").Value = 1.077421197098084
$ws.Range("N6").Value = 1.066470366280739

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.060151775361191
$ws.Range("D7").Value = 1.05906926909791
$ws.Range("E7").Value = 1.065307502674514
$ws.Range("F7").Value = 1.074573830398682
$ws.Range("I7").Value = 1.042718058218711
$ws.Range("J7").Value = 1.064568654323367
$ws.Range("K7").Value = 1.061495586070631
$ws.Range("L7").Value = 1.067718873040512
$ws.Range("M7").Value = 1.076963342660386
$ws.Range("N7").Value = 1.066080463920238

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.058058361176819
$ws.Range("D8").Value = 1.057198544267197
$ws.Range("E8").Value = 1.063388441639612
$ws.Range("F8").Value = 1.072411973460063
$ws.Range("I8").Value = 1.042233883394967
$ws.Range("J8").Value = 1.062935951630019
$ws.Range("K8").Value = 1.059872110039143
$ws.Range("L8").Value = 1.066045556882632
$ws.Range("M8").Value = 1.075045474711102
$ws.Range("N8").Value = 1.064445442601792

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.054351771230171
$ws.Range("D9").Value = 1.053887274060213
$ws.Range("E9").Value = 1.059992853837604
$ws.Range("F9").Value = 1.068590039293027
$ws.Range("I9").Value = 1.04135885138452
$ws.Range("J9").Value = 1.060038464218867
$ws.Range("K9").Value = 1.056992584472888
$ws.Range("L9").Value = 1.063079100268165
$ws.Range("M9").Value = 1.071649845012732
$ws.Range("N9").Value = 1.061543840426169

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.051868139582033
$ws.Range("D10").Value = 1.051669265575735
$ws.Range("E10").Value = 1.057719194054605
$ws.Range("F10").Value = 1.066033077300311
$ws.Range("I10").Value = 1.040760620904636
$ws.Range("J10").Value = 1.058092547750064
$ws.Range("K10").Value = 1.055059798266281
$ws.Range("L10").Value = 1.061088945995441
$ws.Range("M10").Value = 1.069374694812342
$ws.Range("N10").Value = 1.059595160532781

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.050789571466107
$ws.Range("D11").Value = 1.050706234117474
$ws.Range("E11").Value = 1.056732194445032
$ws.Range("F11").Value = 1.064923608562975
$ws.Range("I11").Value = 1.040498010801024
$ws.Range("J11").Value = 1.057246449336696
$ws.Range("K11").Value = 1.054219661920028
$ws.Range("L11").Value = 1.060224104369746
$ws.Range("M11").Value = 1.068386698318233
$ws.Range("N11").Value = 1.058747860562617

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.050388458899488
$ws.Range("D12").Value = 1.050348117414243
$ws.Range("E12").Value = 1.056365194026398
$ws.Range("F12").Value = 1.064511146684357
$ws.Range("I12").Value = 1.040399925349643
$ws.Range("J12").Value = 1.056931634446573
$ws.Range("K12").Value = 1.053907103373692
$ws.Range("L12").Value = 1.059902389260263
$ws.Range("M12").Value = 1.068019274397877
$ws.Range("N12").Value = 1.05843259859924

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.050474521029735
$ws.Range("D13").Value = 1.050424953110722
$ws.Range("E13").Value = 1.056443934449318
$ws.Range("F13").Value = 1.064599637454021
$ws.Range("I13").Value = 1.040420989523893
$ws.Range("J13").Value = 1.056999187735935
$ws.Range("K13").Value = 1.053974170776149
$ws.Range("L13").Value = 1.059971419872473
$ws.Range("M13").Value = 1.068098108077911
$ws.Range("N13").Value = 1.058500247822022

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.050756425330713
$ws.Range("D14").Value = 1.050676640380135
$ws.Range("E14").Value = 1.05670186598859
$ws.Range("F14").Value = 1.064889521664134
$ws.Range("I14").Value = 1.040489914071664
$ws.Range("J14").Value = 1.05722043763682
$ws.Range("K14").Value = 1.054193835876955
$ws.Range("L14").Value = 1.060197521034326
$ws.Range("M14").Value = 1.068356335942904
$ws.Range("N14").Value = 1.058721811923144

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.050930051311634
$ws.Range("D15").Value = 1.050831659373225
$ws.Range("E15").Value = 1.056860734807689
$ws.Range("F15").Value = 1.065068081382245
$ws.Range("I15").Value = 1.040532309046273
$ws.Range("J15").Value = 1.057356685645318
$ws.Range("K15").Value = 1.05432911297877
$ws.Range("L15").Value = 1.060336766257998
$ws.Range("M15").Value = 1.068515380274543
$ws.Range("N15").Value = 1.058858253419441

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.051939653373741
$ws.Range("D16").Value = 1.051733122730577
$ws.Range("E16").Value = 1.057784644594504
$ws.Range("F16").Value = 1.066106659861226
$ws.Range("I16").Value = 1.040777973895556
$ws.Range("J16").Value = 1.058148625855376
$ws.Range("K16").Value = 1.055115486534172
$ws.Range("L16").Value = 1.0611462767021
$ws.Range("M16").Value = 1.069440204071759
$ws.Range("N16").Value = 1.059651318275435

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.052572100680139
$ws.Range("D17").Value = 1.052297878318776
$ws.Range("E17").Value = 1.058363514627079
$ws.Range("F17").Value = 1.066757511783858
$ws.Range("I17").Value = 1.040931113974585
$ws.Range("J17").Value = 1.058644444040284
$ws.Range("K17").Value = 1.055607887101052
$ws.Range("L17").Value = 1.061653226503356
$ws.Range("M17").Value = 1.070019553043422
$ws.Range("N17").Value = 1.060147840579036

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.052940694469918
$ws.Range("D18").Value = 1.052627038700086
$ws.Range("E18").Value = 1.058700920213356
$ws.Range("F18").Value = 1.067136922832492
$ws.Range("I18").Value = 1.041020093561084
$ws.Range("J18").Value = 1.058933309011201
$ws.Range("K18").Value = 1.055894785304175
$ws.Range("L18").Value = 1.061948624118401
$ws.Range("M18").Value = 1.070357204097115
$ws.Range("N18").Value = 1.060437115771345

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.053066324580913
$ws.Range("D19").Value = 1.052739231454661
$ws.Range("E19").Value = 1.058815926496823
$ws.Range("F19").Value = 1.067266255267536
$ws.Range("I19").Value = 1.041050374989068
$ws.Range("J19").Value = 1.059031747571331
$ws.Range("K19").Value = 1.055992557764765
$ws.Range("L19").Value = 1.062049296887378
$ws.Range("M19").Value = 1.070472288342464
$ws.Range("N19").Value = 1.060535694125522

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.052504276401238
$ws.Range("D20").Value = 1.052237311529855
$ws.Range("E20").Value = 1.058301432202238
$ws.Range("F20").Value = 1.0666877043618
$ws.Range("I20").Value = 1.040914719145945
$ws.Range("J20").Value = 1.05859128241038
$ws.Range("K20").Value = 1.055555089388972
$ws.Range("L20").Value = 1.061598866444189
$ws.Range("M20").Value = 1.069957422741214
$ws.Range("N20").Value = 1.06009460345352

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.05067342501338
$ws.Range("D21").Value = 1.050602535987576
$ws.Range("E21").Value = 1.056625922288647
$ws.Range("F21").Value = 1.064804167909134
$ws.Range("I21").Value = 1.040469632456624
$ws.Range("J21").Value = 1.057155299938886
$ws.Range("K21").Value = 1.054129163708443
$ws.Range("L21").Value = 1.060130953054224
$ws.Range("M21").Value = 1.068280306454646
$ws.Range("N21").Value = 1.058656581722208

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.049519487418312
$ws.Range("D22").Value = 1.049572345393043
$ws.Range("E22").Value = 1.055570230873828
$ws.Range("F22").Value = 1.063617849186108
$ws.Range("I22").Value = 1.040186660550062
$ws.Range("J22").Value = 1.056249333126056
$ws.Range("K22").Value = 1.053229761756789
$ws.Range("L22").Value = 1.059205267789281
$ws.Range("M22").Value = 1.067223296965324
$ws.Range("N22").Value = 1.057749328332587

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.050131482076157
$ws.Range("D23").Value = 1.050118694458587
$ws.Range("E23").Value = 1.056130088030617
$ws.Range("F23").Value = 1.064246938972987
$ws.Range("I23").Value = 1.040336967066692
$ws.Range("J23").Value = 1.056729901164836
$ws.Range("K23").Value = 1.053706826650151
$ws.Range("L23").Value = 1.059696255058506
$ws.Range("M23").Value = 1.067783882196327
$ws.Range("N23").Value = 1.058230578833104

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.052534924206692
$ws.Range("D24").Value = 1.052264679833466
$ws.Range("E24").Value = 1.058329485317887
$ws.Range("F24").Value = 1.066719248010465
$ws.Range("I24").Value = 1.040922128327628
$ws.Range("J24").Value = 1.058615304905203
$ws.Range("K24").Value = 1.055578947362586
$ws.Range("L24").Value = 1.06162343033475
$ws.Range("M24").Value = 1.069985497595283
$ws.Range("N24").Value = 1.060118660063041

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.055312181149118
$ws.Range("D25").Value = 1.054745127060202
$ws.Range("E25").Value = 1.060872405022429
$ws.Range("F25").Value = 1.069579642099476
$ws.Range("I25").Value = 1.041587676574825
$ws.Range("J25").Value = 1.060790009983668
$ws.Range("K25").Value = 1.057739283753917
$ws.Range("L25").Value = 1.063848168939282
$ws.Range("M25").Value = 1.072529663265349
$ws.Range("N25").Value = 1.062296453472159

